# Insert two new data rows (2026/01/09 23時台, 2026/01/10 02時台) right after
# the existing 2026/01/09 20時台 row (row 588), pushing all subsequent rows
# down by two. All following rows keep their original values, just shifted.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows starting at row 589 (shifts rows 589..630 down to 591..632)
$insertRange = $ws.Range("A589:D590")
$insertRange.EntireRow.Insert()

# Column A holds dates formatted as plain text (matches the rest of the
# column, which is text, not a real date serial). Force text format so
# Excel doesn't auto-convert the date-like string into a date serial, then
# reset the cell style back to Normal so no stray number format sticks to
# the cell (keeps these cells styleless, like all the other data rows).
$ws.Range("A589:A590").NumberFormat = "@"

# Populate the newly inserted row 589
$ws.Cells.Item(589, 1).Value = "2026/01/09"
$ws.Cells.Item(589, 2).Value = "金"
$ws.Cells.Item(589, 3).Value = 23
$ws.Cells.Item(589, 4).Value = 24

# Populate the newly inserted row 590
$ws.Cells.Item(590, 1).Value = "2026/01/10"
$ws.Cells.Item(590, 2).Value = "土"
$ws.Cells.Item(590, 3).Value = 2
$ws.Cells.Item(590, 4).Value = 24

$ws.Range("A589:A590").Style = "Normal"
